$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D, shifting existing D:K data to F:M
$ws.Columns("D:E").Insert(-4161)

# Re-apply number formatting to the new D:E columns by copying formats
# from F:G (which now hold what used to be in D:E) for each data block
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new D:E columns with the new quarter data
$arr7 = New-Object 'object[,]' 29,2
$arr7[0,0] = 43465
$arr7[0,1] = 43373
$arr7[1,0] = 7600
$arr7[1,1] = 7400
$arr7[2,0] = "NA"
$arr7[2,1] = "NA"
$arr7[3,0] = "NA"
$arr7[3,1] = "NA"
$arr7[4,0] = $null
$arr7[4,1] = $null
$arr7[5,0] = "NA"
$arr7[5,1] = "NA"
$arr7[6,0] = 0
$arr7[6,1] = 0
$arr7[7,0] = 0
$arr7[7,1] = 0
$arr7[8,0] = 0
$arr7[8,1] = 0
$arr7[9,0] = $null
$arr7[9,1] = $null
$arr7[10,0] = 700
$arr7[10,1] = 900
$arr7[11,0] = 7000
$arr7[11,1] = 6500
$arr7[12,0] = $null
$arr7[12,1] = $null
$arr7[13,0] = -4200
$arr7[13,1] = -3500
$arr7[14,0] = 3000
$arr7[14,1] = 3400
$arr7[15,0] = 0
$arr7[15,1] = 0
$arr7[16,0] = 2700
$arr7[16,1] = 3100
$arr7[17,0] = 200
$arr7[17,1] = 500
$arr7[18,0] = 0
$arr7[18,1] = 0
$arr7[19,0] = 2500
$arr7[19,1] = 2500
$arr7[20,0] = 2500
$arr7[20,1] = 2500
$arr7[21,0] = 0
$arr7[21,1] = 0
$arr7[22,0] = 0
$arr7[22,1] = "NA"
$arr7[23,0] = 0
$arr7[23,1] = 0
$arr7[24,0] = 0
$arr7[24,1] = 0
$arr7[25,0] = 4200
$arr7[25,1] = 3500
$arr7[26,0] = 2500
$arr7[26,1] = 2500
$arr7[27,0] = 0
$arr7[27,1] = 0
$arr7[28,0] = 2500
$arr7[28,1] = 2500
$ws.Range("D7:E35").Value = $arr7

$arr38 = New-Object 'object[,]' 40,2
$arr38[0,0] = 43465
$arr38[0,1] = 43373
$arr38[1,0] = $null
$arr38[1,1] = $null
$arr38[2,0] = $null
$arr38[2,1] = $null
$arr38[3,0] = 17400
$arr38[3,1] = 22300
$arr38[4,0] = 33400
$arr38[4,1] = 38800
$arr38[5,0] = 0
$arr38[5,1] = 0
$arr38[6,0] = 0
$arr38[6,1] = 0
$arr38[7,0] = 0
$arr38[7,1] = 0
$arr38[8,0] = 0
$arr38[8,1] = 0
$arr38[9,0] = 0
$arr38[9,1] = 0
$arr38[10,0] = 14400
$arr38[10,1] = 14400
$arr38[11,0] = 7400
$arr38[11,1] = 7500
$arr38[12,0] = 0
$arr38[12,1] = 0
$arr38[13,0] = 0
$arr38[13,1] = 0
$arr38[14,0] = 0
$arr38[14,1] = 0
$arr38[15,0] = 0
$arr38[15,1] = 0
$arr38[16,0] = 794200
$arr38[16,1] = 790500
$arr38[17,0] = $null
$arr38[17,1] = $null
$arr38[18,0] = $null
$arr38[18,1] = $null
$arr38[19,0] = 200
$arr38[19,1] = 100
$arr38[20,0] = 0
$arr38[20,1] = 0
$arr38[21,0] = 6400
$arr38[21,1] = 7000
$arr38[22,0] = 0
$arr38[22,1] = 0
$arr38[23,0] = 0
$arr38[23,1] = 0
$arr38[24,0] = 0
$arr38[24,1] = 0
$arr38[25,0] = 0
$arr38[25,1] = 0
$arr38[26,0] = 0
$arr38[26,1] = 0
$arr38[27,0] = 0
$arr38[27,1] = 0
$arr38[28,0] = 708300
$arr38[28,1] = 709300
$arr38[29,0] = $null
$arr38[29,1] = $null
$arr38[30,0] = 0
$arr38[30,1] = 0
$arr38[31,0] = 0
$arr38[31,1] = 0
$arr38[32,0] = 0
$arr38[32,1] = 0
$arr38[33,0] = 0
$arr38[33,1] = 0
$arr38[34,0] = 57400
$arr38[34,1] = 55600
$arr38[35,0] = 0
$arr38[35,1] = 0
$arr38[36,0] = 0
$arr38[36,1] = 0
$arr38[37,0] = 0
$arr38[37,1] = 0
$arr38[38,0] = 85800
$arr38[38,1] = 81200
$arr38[39,0] = 0
$arr38[39,1] = 0
$ws.Range("D38:E77").Value = $arr38

$arr80 = New-Object 'object[,]' 23,2
$arr80[0,0] = 43465
$arr80[0,1] = 43373
$arr80[1,0] = 2500
$arr80[1,1] = 2500
$arr80[2,0] = $null
$arr80[2,1] = $null
$arr80[3,0] = 300
$arr80[3,1] = 300
$arr80[4,0] = 0
$arr80[4,1] = 0
$arr80[5,0] = 0
$arr80[5,1] = 0
$arr80[6,0] = 0
$arr80[6,1] = 0
$arr80[7,0] = 0
$arr80[7,1] = 0
$arr80[8,0] = 0
$arr80[8,1] = 0
$arr80[9,0] = 2900
$arr80[9,1] = 4100
$arr80[10,0] = $null
$arr80[10,1] = $null
$arr80[11,0] = -200
$arr80[11,1] = -100
$arr80[12,0] = 0
$arr80[12,1] = 0
$arr80[13,0] = 0
$arr80[13,1] = 0
$arr80[14,0] = -10700
$arr80[14,1] = -15300
$arr80[15,0] = $null
$arr80[15,1] = $null
$arr80[16,0] = -800
$arr80[16,1] = -800
$arr80[17,0] = 0
$arr80[17,1] = 0
$arr80[18,0] = 0
$arr80[18,1] = 0
$arr80[19,0] = 0
$arr80[19,1] = 0
$arr80[20,0] = -1200
$arr80[20,1] = -600
$arr80[21,0] = 0
$arr80[21,1] = 0
$arr80[22,0] = -9000
$arr80[22,1] = -11700
$ws.Range("D80:E102").Value = $arr80
